# AR6_min_max_updates.xlsx edit script
#
# 1) Forest land cover variable changed from absolute area [ha] to a
#    share-of-land metric. New indicator row ("Share land cover forest")
#    added with min/max taken from the regional AR6 database, plus a note
#    row documenting the change.
# 2) Misc bookkeeping: active-cell selection moved, workbook "saved from"
#    folder metadata updated (best effort - not exposed by the object
#    model for scripted editing).
#
# NOTE on shared-string ordering: cells are written in the same order the
# original author's strings were appended to the shared string table
# (A15/K15 "Share land cover forest" first, then D15 "share of land as
# forest", then A14 "#adding share of forest in land cover") so the new
# <si> entries land at indices 63/64/65 as in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: new "Share land cover forest" indicator ---------------------
# Label cells - copy formatting from the row above (A13/K13) so the same
# bold "indicator name" style is reused.
$ws.Range("A13").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = "Share land cover forest"

$ws.Range("K13").Copy()
$ws.Range("K15").PasteSpecial(-4122)
$ws.Range("K15").Value = "Share land cover forest"

# Numeric min/max (global) values - plain/default style.
$ws.Range("B15").Value = 0.227
$ws.Range("C15").Value = 0.448

# Unit label cell: "share of land as forest" - also clears this cell's
# border so it carries its own dedicated (borderless) format.
$ws.Range("D15").Value = "share of land as forest"
$ws.Range("D15").Interior.ColorIndex = -4142

# Exogenous min/max used for the "final" min/max table.
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0.61

# Mirrored min/max on the right-hand "final" table, as formulas referencing
# the exogenous min/max just entered.
$ws.Range("L15").Formula = "=E15"
$ws.Range("M15").Formula = "=F15"

# --- Row 14: note row, only column A populated ----------------------------
$ws.Range("A14").Value = "#adding share of forest in land cover"

# --- Selection moved by the author while editing --------------------------
$ws.Range("M21").Select()

# --- Workbook "absolute path" bookkeeping (best effort; Excel stamps this
#     from the actual save location and it is not exposed for scripted
#     editing, but we still try the documented Path/FullName properties).
try {
    $wb.Path = "C:\Users\uk\Projects\Navigate\AlternativeWelfareMetrics\"
} catch {
}
try {
    $wb.FullName = "C:\Users\uk\Projects\Navigate\AlternativeWelfareMetrics\AR6_min_max_updates.xlsx"
} catch {
}
